# Swap the species-record data between row 14 and row 15 on the "Artfynd"
# sheet, as described by the source diff. The two rows share the same
# "context" columns (county, municipality, dates, reporter, etc.) which
# stay untouched; only the record-specific columns are exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Capture current (pre-edit) values for row 14 ---------------------
$a14 = $ws.Range("A14").Value()
$b14 = $ws.Range("B14").Value()
$c14 = $ws.Range("C14").Value()
$d14 = $ws.Range("D14").Value()
$e14 = $ws.Range("E14").Value()
$f14 = $ws.Range("F14").Value()
$g14 = $ws.Range("G14").Value()
$h14 = $ws.Range("H14").Value()
$i14 = $ws.Range("I14").Value()
$j14 = $ws.Range("J14").Value()
$k14 = $ws.Range("K14").Value()
$m14 = $ws.Range("M14").Value()
$p14 = $ws.Range("P14").Value()
$q14 = $ws.Range("Q14").Value()
$r14 = $ws.Range("R14").Value()
$s14 = $ws.Range("S14").Value()

# --- Capture current (pre-edit) values for row 15 ---------------------
$a15 = $ws.Range("A15").Value()
$b15 = $ws.Range("B15").Value()
$c15 = $ws.Range("C15").Value()
$d15 = $ws.Range("D15").Value()
$e15 = $ws.Range("E15").Value()
$f15 = $ws.Range("F15").Value()
$g15 = $ws.Range("G15").Value()
$h15 = $ws.Range("H15").Value()
$i15 = $ws.Range("I15").Value()
$j15 = $ws.Range("J15").Value()
$k15 = $ws.Range("K15").Value()
$m15 = $ws.Range("M15").Value()
$p15 = $ws.Range("P15").Value()
$q15 = $ws.Range("Q15").Value()
$r15 = $ws.Range("R15").Value()
$s15 = $ws.Range("S15").Value()

# --- Write row 15's old data into row 14 -------------------------------
$ws.Range("A14").Value = $a15
$ws.Range("B14").Value = $b15
$ws.Range("C14").Value = $c15
$ws.Range("D14").Value = $d15
$ws.Range("E14").Value = $e15
$ws.Range("F14").Value = $f15
$ws.Range("G14").Value = $g15
$ws.Range("H14").Value = $h15
$ws.Range("I14").Value = $i15
$ws.Range("J14").Value = $j15
$ws.Range("K14").Value = $k15
$ws.Range("M14").ClearContents()
$ws.Range("P14").Value = $p15
$ws.Range("Q14").Value = $q15
$ws.Range("R14").Value = $r15
$ws.Range("S14").Value = $s15

# --- Write row 14's old data into row 15 -------------------------------
$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("C15").Value = $c14
$ws.Range("D15").Value = $d14
$ws.Range("E15").Value = $e14
$ws.Range("F15").Value = $f14
$ws.Range("G15").Value = $g14
$ws.Range("H15").Value = $h14
$ws.Range("I15").Value = $i14
$ws.Range("J15").ClearContents()
$ws.Range("K15").Value = $k14
$ws.Range("M15").Value = $m14
$ws.Range("P15").Value = $p14
$ws.Range("Q15").Value = $q14
$ws.Range("R15").Value = $r14
$ws.Range("S15").Value = $s14
